# Update countries & provincias Spain
#
# Refreshes the "Pais" sheet with the newer data snapshot:
#  - Updates the "last updated" timestamp in A1.
#  - Updates the Covid case counters (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Muertes hoy, Muertes) for the countries
#    whose figures moved since the previous snapshot.
#  - A few countries overtook their neighbours in the "Casos totales"
#    ranking (Sudafrica/Alemania, Nigeria/Rumania/Armenia,
#    Bulgaria/Etiopia, Islas Malvinas/Groenlandia), so those rows swap
#    places/names to keep the table sorted by total cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Julio de 2020 a las 00:49"

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 3026612
$ws.Cells.Item(4, 3).Value = 43684
$ws.Cells.Item(4, 4).Value = 1308563
$ws.Cells.Item(4, 5).Value = 1585213
$ws.Cells.Item(4, 7).Value = 267
$ws.Cells.Item(4, 8).Value = 132836

# Brasil (row 5)
$ws.Cells.Item(5, 2).Value = 1623284
$ws.Cells.Item(5, 3).Value = 18699
$ws.Cells.Item(5, 5).Value = 579182
$ws.Cells.Item(5, 7).Value = 587
$ws.Cells.Item(5, 8).Value = 65487

# Row 18: Alemania -> Sudafrica (Sudafrica overtakes Alemania)
$ws.Cells.Item(18, 1).Value = "Sudafrica"
$ws.Cells.Item(18, 2).Value = 205721
$ws.Cells.Item(18, 3).Value = 8971
$ws.Cells.Item(18, 4).Value = 97848
$ws.Cells.Item(18, 5).Value = 104563
$ws.Cells.Item(18, 7).Value = 111
$ws.Cells.Item(18, 8).Value = 3310

# Row 19: Sudafrica -> Alemania
$ws.Cells.Item(19, 1).Value = "Alemania"
$ws.Cells.Item(19, 2).Value = 198057
$ws.Cells.Item(19, 3).Value = 499
$ws.Cells.Item(19, 4).Value = 182200
$ws.Cells.Item(19, 5).Value = 6765
$ws.Cells.Item(19, 7).Value = 6
$ws.Cells.Item(19, 8).Value = 9092

# Colombia (row 22)
$ws.Cells.Item(22, 2).Value = 120281
$ws.Cells.Item(22, 3).Value = 3171
$ws.Cells.Item(22, 4).Value = 50370
$ws.Cells.Item(22, 5).Value = 65701
$ws.Cells.Item(22, 7).Value = 146
$ws.Cells.Item(22, 8).Value = 4210

# Row 51: Rumania -> Nigeria (Nigeria overtakes Rumania & Armenia)
$ws.Cells.Item(51, 1).Value = "Nigeria"
$ws.Cells.Item(51, 2).Value = 29286
$ws.Cells.Item(51, 3).Value = 575
$ws.Cells.Item(51, 4).Value = 11828
$ws.Cells.Item(51, 5).Value = 16804
$ws.Cells.Item(51, 7).Value = 9
$ws.Cells.Item(51, 8).Value = 654

# Row 52: Armenia -> Rumania
$ws.Cells.Item(52, 1).Value = "Rumania"
$ws.Cells.Item(52, 2).Value = 29223
$ws.Cells.Item(52, 3).Value = 250
$ws.Cells.Item(52, 4).Value = 20213
$ws.Cells.Item(52, 5).Value = 7242
$ws.Cells.Item(52, 7).Value = 18
$ws.Cells.Item(52, 8).Value = 1768

# Row 53: Nigeria -> Armenia
$ws.Cells.Item(53, 1).Value = "Armenia"
$ws.Cells.Item(53, 2).Value = 28936
$ws.Cells.Item(53, 3).Value = 330
$ws.Cells.Item(53, 4).Value = 16302
$ws.Cells.Item(53, 5).Value = 12143
$ws.Cells.Item(53, 7).Value = 7
$ws.Cells.Item(53, 8).Value = 491

# Barein-area country (row 59)
$ws.Cells.Item(59, 2).Value = 19775
$ws.Cells.Item(59, 3).Value = 253
$ws.Cells.Item(59, 4).Value = 17124
$ws.Cells.Item(59, 5).Value = 1674

# row 69
$ws.Cells.Item(69, 2).Value = 12566
$ws.Cells.Item(69, 3).Value = 51
$ws.Cells.Item(69, 5).Value = 4343

# row 72
$ws.Cells.Item(72, 2).Value = 9894
$ws.Cells.Item(72, 3).Value = 127
$ws.Cells.Item(72, 4).Value = 4899
$ws.Cells.Item(72, 5).Value = 4379
$ws.Cells.Item(72, 7).Value = 8
$ws.Cells.Item(72, 8).Value = 616

# Row 86: Etiopia -> Bulgaria (Bulgaria overtakes Etiopia)
$ws.Cells.Item(86, 1).Value = "Bulgaria"
$ws.Cells.Item(86, 2).Value = 5914
$ws.Cells.Item(86, 3).Value = 174
$ws.Cells.Item(86, 4).Value = 3000
$ws.Cells.Item(86, 5).Value = 2664
$ws.Cells.Item(86, 7).Value = 4
$ws.Cells.Item(86, 8).Value = 250

# Row 87: Bulgaria -> Etiopia
$ws.Cells.Item(87, 1).Value = "Etiopia"
$ws.Cells.Item(87, 2).Value = 5846
$ws.Cells.Item(87, 4).Value = 2430
$ws.Cells.Item(87, 5).Value = 3313
$ws.Cells.Item(87, 8).Value = 103

# row 88
$ws.Cells.Item(88, 2).Value = 5743
$ws.Cells.Item(88, 3).Value = 123
$ws.Cells.Item(88, 4).Value = 2574
$ws.Cells.Item(88, 5).Value = 3123
$ws.Cells.Item(88, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 46

# row 189
$ws.Cells.Item(189, 2).Value = 70
$ws.Cells.Item(189, 3).Value = 2
$ws.Cells.Item(189, 5).Value = 44

# Row 209: Groenlandia -> Islas Malvinas (Islas Malvinas overtakes Groenlandia)
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"

# Row 210: Islas Malvinas -> Groenlandia
$ws.Cells.Item(210, 1).Value = "Groenlandia"
